$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 77.011819
$ws.Range("H2").Value = 231.035457
$ws.Range("I2").Value = 0.363766179891216
$ws.Range("J2").Value = 0.3637661798912161
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1451143333333333
$ws.Range("N2").Value = 0.435343
$ws.Range("O2").Value = 0.140827770705994
$ws.Range("P2").Value = 0.1408277707059941
$ws.Range("Q2").Value = 11.17551877297233
$ws.Range("R2").Value = 100.579668956751
$ws.Range("S2").Value = 0.05122838017231555
$ws.Range("T2").Value = 0.05122838017231557
$ws.Range("G3").Value = 77.011819
$ws.Range("H3").Value = 231.035457
$ws.Range("I3").Value = 0.363766179891216
$ws.Range("J3").Value = 0.3637661798912161
$ws.Range("O3").Value = 0.0422072807203407
$ws.Range("P3").Value = 0.0422072807203407
$ws.Range("Q3").Value = 3.349398031948
$ws.Range("R3").Value = 30.144582287532
$ws.Range("S3").Value = 0.01535358127123451
$ws.Range("T3").Value = 0.01535358127123451
$ws.Range("G4").Value = 77.011819
$ws.Range("H4").Value = 231.035457
$ws.Range("I4").Value = 0.363766179891216
$ws.Range("J4").Value = 0.3637661798912161
$ws.Range("M4").Value = 0.841832
$ws.Range("N4").Value = 2.525496
$ws.Range("O4").Value = 0.8169649485736653
$ws.Range("P4").Value = 0.8169649485736653
$ws.Range("Q4").Value = 64.831013612408
$ws.Range("R4").Value = 583.479122511672
$ws.Range("S4").Value = 0.297184218447666
$ws.Range("T4").Value = 0.2971842184476661
$ws.Range("I5").Value = 0.5327430638656547
$ws.Range("J5").Value = 0.532743063865655
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1451143333333333
$ws.Range("N5").Value = 0.435343
$ws.Range("O5").Value = 0.140827770705994
$ws.Range("P5").Value = 0.1408277707059941
$ws.Range("Q5").Value = 16.36677745353311
$ws.Range("R5").Value = 147.300997081798
$ws.Range("S5").Value = 0.07502501804328116
$ws.Range("T5").Value = 0.07502501804328121
$ws.Range("I6").Value = 0.5327430638656547
$ws.Range("J6").Value = 0.532743063865655
$ws.Range("O6").Value = 0.0422072807203407
$ws.Range("P6").Value = 0.0422072807203407
$ws.Range("S6").Value = 0.02248563604839208
$ws.Range("T6").Value = 0.0224856360483921
$ws.Range("I7").Value = 0.5327430638656547
$ws.Range("J7").Value = 0.532743063865655
$ws.Range("M7").Value = 0.841832
$ws.Range("N7").Value = 2.525496
$ws.Range("O7").Value = 0.8169649485736653
$ws.Range("P7").Value = 0.8169649485736653
$ws.Range("Q7").Value = 94.94635492425067
$ws.Range("R7").Value = 854.517194318256
$ws.Range("S7").Value = 0.4352324097739815
$ws.Range("T7").Value = 0.4352324097739817
$ws.Range("G8").Value = 0.080633
$ws.Range("H8").Value = 0.241899
$ws.Range("I8").Value = 0.0003808708684464188
$ws.Range("J8").Value = 0.0003808708684464189
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1451143333333333
$ws.Range("N8").Value = 0.435343
$ws.Range("O8").Value = 0.140827770705994
$ws.Range("P8").Value = 0.1408277707059941
$ws.Range("Q8").Value = 0.01170100403966667
$ws.Range("R8").Value = 0.105309036357
$ws.Range("S8").Value = 0.00005363719533016509
$ws.Range("T8").Value = 0.00005363719533016512
$ws.Range("G9").Value = 0.080633
$ws.Range("H9").Value = 0.241899
$ws.Range("I9").Value = 0.0003808708684464188
$ws.Range("J9").Value = 0.0003808708684464189
$ws.Range("O9").Value = 0.0422072807203407
$ws.Range("P9").Value = 0.0422072807203407
$ws.Range("Q9").Value = 0.003506890436
$ws.Range("R9").Value = 0.031562013924
$ws.Range("S9").Value = 0.00001607552366271795
$ws.Range("T9").Value = 0.00001607552366271796
$ws.Range("G10").Value = 0.080633
$ws.Range("H10").Value = 0.241899
$ws.Range("I10").Value = 0.0003808708684464188
$ws.Range("J10").Value = 0.0003808708684464189
$ws.Range("M10").Value = 0.841832
$ws.Range("N10").Value = 2.525496
$ws.Range("O10").Value = 0.8169649485736653
$ws.Range("P10").Value = 0.8169649485736653
$ws.Range("Q10").Value = 0.067879439656
$ws.Range("R10").Value = 0.610914956904
$ws.Range("S10").Value = 0.0003111581494535358
$ws.Range("T10").Value = 0.0003111581494535359
$ws.Range("G11").Value = 21.34689033333333
$ws.Range("H11").Value = 64.040671
$ws.Range("I11").Value = 0.1008322728893521
$ws.Range("J11").Value = 0.1008322728893522
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1451143333333333
$ws.Range("N11").Value = 0.435343
$ws.Range("O11").Value = 0.140827770705994
$ws.Range("P11").Value = 0.1408277707059941
$ws.Range("Q11").Value = 3.097739759461444
$ws.Range("R11").Value = 27.879657835153
$ws.Range("S11").Value = 0.0141999842062259
$ws.Range("T11").Value = 0.01419998420622591
$ws.Range("G12").Value = 21.34689033333333
$ws.Range("H12").Value = 64.040671
$ws.Range("I12").Value = 0.1008322728893521
$ws.Range("J12").Value = 0.1008322728893522
$ws.Range("O12").Value = 0.0422072807203407
$ws.Range("P12").Value = 0.0422072807203407
$ws.Range("Q12").Value = 0.9284189543773335
$ws.Range("R12").Value = 8.355770589396
$ws.Range("S12").Value = 0.004255856047510884
$ws.Range("T12").Value = 0.004255856047510886
$ws.Range("G13").Value = 21.34689033333333
$ws.Range("H13").Value = 64.040671
$ws.Range("I13").Value = 0.1008322728893521
$ws.Range("J13").Value = 0.1008322728893522
$ws.Range("M13").Value = 0.841832
$ws.Range("N13").Value = 2.525496
$ws.Range("O13").Value = 0.8169649485736653
$ws.Range("P13").Value = 0.8169649485736653
$ws.Range("Q13").Value = 17.97049538309067
$ws.Range("R13").Value = 161.734458447816
$ws.Range("S13").Value = 0.08237643263561535
$ws.Range("T13").Value = 0.08237643263561538
$ws.Range("G14").Value = 0.465901
$ws.Range("H14").Value = 1.397703
$ws.Range("I14").Value = 0.002200688532983456
$ws.Range("J14").Value = 0.002200688532983456
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1451143333333333
$ws.Range("N14").Value = 0.435343
$ws.Range("O14").Value = 0.140827770705994
$ws.Range("P14").Value = 0.1408277707059941
$ws.Range("Q14").Value = 0.06760891301433332
$ws.Range("R14").Value = 0.6084802171289999
$ws.Range("S14").Value = 0.0003099180601183045
$ws.Range("T14").Value = 0.0003099180601183046
$ws.Range("G15").Value = 0.465901
$ws.Range("H15").Value = 1.397703
$ws.Range("I15").Value = 0.002200688532983456
$ws.Range("J15").Value = 0.002200688532983456
$ws.Range("O15").Value = 0.0422072807203407
$ws.Range("P15").Value = 0.0422072807203407
$ws.Range("Q15").Value = 0.020262966292
$ws.Range("R15").Value = 0.182366696628
$ws.Range("S15").Value = 0.00009288507868966746
$ws.Range("T15").Value = 0.0000928850786896675
$ws.Range("G16").Value = 0.465901
$ws.Range("H16").Value = 1.397703
$ws.Range("I16").Value = 0.002200688532983456
$ws.Range("J16").Value = 0.002200688532983456
$ws.Range("M16").Value = 0.841832
$ws.Range("N16").Value = 2.525496
$ws.Range("O16").Value = 0.8169649485736653
$ws.Range("P16").Value = 0.8169649485736653
$ws.Range("Q16").Value = 0.392210370632
$ws.Range("R16").Value = 3.529893335688
$ws.Range("S16").Value = 0.001797885394175484
$ws.Range("T16").Value = 0.001797885394175484
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.01628533333333333
$ws.Range("H17").Value = 0.048856
$ws.Range("I17").Value = 0.00007692395234712933
$ws.Range("J17").Value = 0.00007692395234712935
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1451143333333333
$ws.Range("N17").Value = 0.435343
$ws.Range("O17").Value = 0.140827770705994
$ws.Range("P17").Value = 0.1408277707059941
$ws.Range("Q17").Value = 0.002363235289777777
$ws.Range("R17").Value = 0.021269117608
$ws.Range("S17").Value = 0.00001083302872294034
$ws.Range("T17").Value = 0.00001083302872294034
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.01628533333333333
$ws.Range("H18").Value = 0.048856
$ws.Range("I18").Value = 0.00007692395234712933
$ws.Range("J18").Value = 0.00007692395234712935
$ws.Range("O18").Value = 0.0422072807203407
$ws.Range("P18").Value = 0.0422072807203407
$ws.Range("Q18").Value = 0.0007082817173333333
$ws.Range("R18").Value = 0.006374535456
$ws.Range("S18").Value = 0.000003246750850833398
$ws.Range("T18").Value = 0.0000032467508508334
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.01628533333333333
$ws.Range("H19").Value = 0.048856
$ws.Range("I19").Value = 0.00007692395234712933
$ws.Range("J19").Value = 0.00007692395234712935
$ws.Range("M19").Value = 0.841832
$ws.Range("N19").Value = 2.525496
$ws.Range("O19").Value = 0.8169649485736653
$ws.Range("P19").Value = 0.8169649485736653
$ws.Range("Q19").Value = 0.01370951473066667
$ws.Range("R19").Value = 0.123385632576
$ws.Range("S19").Value = 0.0000628441727733556
$ws.Range("T19").Value = 0.00006284417277335561

Write-Host "applied updates"
